# The workbook's "locations" sheet has two data rows:
#   Row 2: "CTY ABC" sample record
#   Row 3: "Ho kinh doanh Tran Van A" demo record
#
# This edit removes the "CTY ABC" row entirely (row 2), which shifts the
# demo record up to become the new row 2, and then refreshes that
# surviving record's id / code / createdAt to newly generated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the first data row (the "CTY ABC" record); this shifts row 3 up to row 2.
$ws.Rows(2).Delete()

# Update the now-surviving row 2 (formerly row 3) with refreshed id/code/createdAt.
$ws.Range("A2").Value = "zdxwxqmyo3o"
$ws.Range("B2").Value = "DEMOT9QZ"
$ws.Range("F2").Value = "2025-08-13T00:50:53.687Z"
